$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1511
$ws.Range("F3").Value = 841
$ws.Range("F4").Value = 436
$ws.Range("F5").Value = 886
$ws.Range("F6").Value = 497
$ws.Range("F7").Value = 7463
$ws.Range("F11").Value = 5437
$ws.Range("F12").Value = 561
$ws.Range("F13").Value = 241
$ws.Range("F15").Value = 7502
$ws.Range("F16").Value = 8825
$ws.Range("F17").Value = 1136
$ws.Range("F18").Value = 877
$ws.Range("F19").Value = 4404
$ws.Range("F20").Value = 661
$ws.Range("F21").Value = 206
$ws.Range("F23").Value = 281
$ws.Range("F25").Value = 1184
$ws.Range("F26").Value = 96
$ws.Range("F27").Value = 1639
$ws.Range("F28").Value = 698
$ws.Range("F29").Value = 896
$ws.Range("F30").Value = 1
$ws.Range("F31").Value = 1858
$ws.Range("F32").Value = 325
$ws.Range("F33").Value = 2255
$ws.Range("F35").Value = 105
$ws.Range("F36").Value = 1432
$ws.Range("F39").Value = 787
$ws.Range("F40").Value = 389
$ws.Range("F41").Value = 4043
$ws.Range("F42").Value = 187
$ws.Range("F44").Value = 411
$ws.Range("F46").Value = 13
$ws.Range("F48").Value = 164
$ws.Range("F49").Value = 4062

# 演出 (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 46
$ws.Range("F25").Value = 61

# 本地生活 (sheet3)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5160

# 全部类型 (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1511
$ws.Range("F5").Value = 841
$ws.Range("F6").Value = 886
$ws.Range("F11").Value = 5437
$ws.Range("F12").Value = 561
$ws.Range("F13").Value = 7502
$ws.Range("F15").Value = 1136
$ws.Range("F16").Value = 877
$ws.Range("F17").Value = 4404
$ws.Range("F18").Value = 661
$ws.Range("F19").Value = 206
$ws.Range("F21").Value = 281
$ws.Range("F24").Value = 46
$ws.Range("F25").Value = 1184
$ws.Range("F26").Value = 96
$ws.Range("F27").Value = 1639
$ws.Range("F28").Value = 1858
$ws.Range("F29").Value = 325
$ws.Range("F30").Value = 2255
$ws.Range("F37").Value = 787
$ws.Range("F39").Value = 61
$ws.Range("F40").Value = 389
$ws.Range("F41").Value = 4043
$ws.Range("F43").Value = 187
$ws.Range("F45").Value = 411
$ws.Range("F48").Value = 164
$ws.Range("F49").Value = 4062
